# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K" in row 1) holds newly (re)computed
# values ("s_vals") for each data row (rows 2-40). Write the recalculated
# values in place, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (the "K" column), for rows 2 through 40 in order.
$newKValues = @(2,1,2,2,0,1,1,1,1,2,0,1,2,2,2,2,2,1,1,1,1,0,1,1,1,0,2,3,1,1,0,2,1,1,0,1,1,1,1)

$row = 2
foreach ($val in $newKValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
